$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Row 3: clear the "Application-JS-4.0" subject/build label in C3.
$ws.Range("C3").Value = ""

# Row 3: add a ScrollPage step before the final ClickRunTest/validate4 in the
# VT200_0004 automation script (column G).
$ws.Range("G3").Value = "wait(3);`nvalidate1;`nlink_Click(Application_test_link);`nvalidate2;`nSelectTestToRun(VT200_0004_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nScrollPage(runtest_bottom_xpath);`nClickRunTest(runtest_bottom_xpath);`nvalidate4;"

# Row 4: same ScrollPage addition for the VT200_0005 automation script (column G).
$ws.Range("G4").Value = "wait(3);`nvalidate1;`nlink_Click(Application_test_link);`nvalidate2;`nSelectTestToRun(VT200_0005_string);`nClickRunTest(runtest_top_xpath);`nvalidate3;`nScrollPage(runtest_bottom_xpath);`nClickRunTest(runtest_bottom_xpath);`nvalidate4;"

# Row 2: update the start-page URL embedded in G2's automation script text.
$ws.Range("G2").Value = "wait(3);`nSetStartPage(http://127.0.0.1:8082/app/);"

# Update the sheet selection to J2:J12 (active cell J2).
$ws.Range("J2:J12").Select()
